# Update column C ("Förändrad") date values for rows 2 through 91
# from Excel serial date 45204 (2023-10-05) to 45205 (2023-10-06).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 91; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45204) {
        $cell.Value2 = 45205
    }
}
